# The deck ships two theme parts:
#   ppt/theme/theme1.xml  -> used by the (only) slide master ("Integral" theme)
#   ppt/theme/theme2.xml  -> used by the notes master ("Office Theme")
#
# The target revision swaps their contents: the slide-master theme becomes the
# stock "Office Theme" colours, and the notes-master theme becomes the
# "Integral" colours. Font scheme / format scheme (fills, lines, effects) are
# byte-identical between the two themes already, so only the 12-slot colour
# scheme (and the cosmetic theme/colour-scheme names) actually differ.
#
# Apply the reachable part of that swap through the real PowerPoint object
# model: each slide's ThemeColorScheme maps straight onto the slide master's
# theme (theme1.xml), so push the "Office Theme" palette onto it, in the
# standard clrScheme slot order (dk1, lt1, dk2, lt2, accent1-6, hlink,
# folHlink).

$p = $ppt.ActivePresentation

$officeThemeColors = @(
    0,         # 1  dk1      000000
    16777215,  # 2  lt1      FFFFFF
    6968388,   # 3  dk2      44546A
    15132391,  # 4  lt2      E7E6E6
    13998939,  # 5  accent1  5B9BD5
    3243501,   # 6  accent2  ED7D31
    10855845,  # 7  accent3  A5A5A5
    49407,     # 8  accent4  FFC000
    12874308,  # 9  accent5  4472C4
    4697456,   # 10 accent6  70AD47
    12673797,  # 11 hlink    0563C1
    7491477    # 12 folHlink 954F72
)

$slide = $p.Slides.Item(1)
$themeColors = $slide.ThemeColorScheme

for ($i = 1; $i -le 12; $i++) {
    $themeColors.Item($i).RGB = $officeThemeColors[$i - 1]
}

# Best-effort: also try to rename the design/theme to match "Office Theme" /
# "Office" so the scheme metadata lines up with the swapped palette (no-ops
# harmlessly if the host does not support renaming designs).
try {
    $design = $p.Designs.Item(1)
    $design.Name = "Office Theme"
} catch {
}
